$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2166.5   # was 1944.3334
$ws.Range("I62").Value = 1899.75   # was 1728.4286
$ws.Range("K62").Value = 1899.75   # was 1728.4286
$ws.Range("M62").Value = -1275.75   # was -1104.4286
$ws.Range("H65").Value = 2166.5   # was 1944.3334
$ws.Range("I65").Value = 1899.75   # was 1728.4286
$ws.Range("K65").Value = 9498.75   # was 8642.143
$ws.Range("M65").Value = -6378.75   # was -5522.143
$ws.Range("H96").Value = 469.77777   # was 347.22223
$ws.Range("I96").Value = 473   # was 219.2
$ws.Range("J96").Value = 468.16666   # was 507.25
$ws.Range("K96").Value = 1419   # was 657.5999999999999
$ws.Range("L96").Value = 1404.49998   # was 1521.75
$ws.Range("M96").Value = -46   # was 715.4000000000001
$ws.Range("N96").Value = -4150.499980000001   # was -4267.75
$ws.Range("H100").Value = 18520758   # was 18519940
$ws.Range("I100").Value = 27779636   # was 23810780
$ws.Range("J100").Value = 2999.6667   # was 2000
$ws.Range("K100").Value = 27779636   # was 23810780
$ws.Range("L100").Value = 2999.6667   # was 2000
$ws.Range("M100").Value = -27779095   # was -23810239
$ws.Range("N100").Value = -4081.6667   # was -3082
$ws.Range("H116").Value = 5960.9614   # was 8100.3125
$ws.Range("I116").Value = 9538.846   # was 11640.5
$ws.Range("J116").Value = 2383.077   # was 2200
$ws.Range("K116").Value = 9538.846   # was 11640.5
$ws.Range("L116").Value = 2383.077   # was 2200
$ws.Range("M116").Value = -6096.846   # was -8198.5
$ws.Range("N116").Value = -9267.077000000001   # was -9084
$ws.Range("H137").Value = 1911.1177   # was 1772
$ws.Range("I137").Value = 1185.4667   # was 1260.1538
$ws.Range("J137").Value = 2484   # was 2027.9231
$ws.Range("K137").Value = 3556.4001   # was 3780.4614
$ws.Range("L137").Value = 7452   # was 6083.7693
$ws.Range("M137").Value = -1006.4001   # was -1230.4614
$ws.Range("N137").Value = -12552   # was -11183.7693
$ws.Range("H138").Value = 2782.1875   # was 2767.7334
$ws.Range("J138").Value = 4005.2   # was 4130.975
$ws.Range("L138").Value = 12015.6   # was 12392.925
$ws.Range("N138").Value = -22295.6   # was -22672.925

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9470.944   # was 6881.077
$ws.Range("I45").Value = 11723.357   # was 8787.315000000001
$ws.Range("J45").Value = 1587.5   # was 1707
$ws.Range("K45").Value = 11723.357   # was 8787.315000000001
$ws.Range("L45").Value = 1587.5   # was 1707
$ws.Range("M45").Value = -11346.357   # was -8410.315000000001
$ws.Range("N45").Value = -2341.5   # was -2461
$ws.Range("H88").Value = 1000000000   # was 125002216
$ws.Range("I88").Value = 0   # was 2526
$ws.Range("J88").Value = 1000000000   # was 250001900
$ws.Range("K88").Value = 0   # was 2526
$ws.Range("L88").Value = 1000000000   # was 250001900
$ws.Range("M88").Value = $null   # was -2120
$ws.Range("N88").Value = -1000000812   # was -250002712
$ws.Range("H91").Value = 1000000000   # was 125002216
$ws.Range("I91").Value = 0   # was 2526
$ws.Range("J91").Value = 1000000000   # was 250001900
$ws.Range("K91").Value = 0   # was 2526
$ws.Range("L91").Value = 1000000000   # was 250001900
$ws.Range("M91").Value = $null   # was -1122
$ws.Range("N91").Value = -1000002808   # was -250004708
$ws.Range("H122").Value = 3668621.2   # was 4279825
$ws.Range("I122").Value = 3668621.2   # was 4279825
$ws.Range("K122").Value = 11005863.6   # was 12839475
$ws.Range("M122").Value = -11003413.6   # was -12837025
$ws.Range("H132").Value = 8474.75   # was 6568.4287
$ws.Range("J132").Value = 9999.5   # was 6415.8
$ws.Range("L132").Value = 29998.5   # was 19247.4
$ws.Range("N132").Value = -35058.5   # was -24307.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2393.5334   # was 2425.0625
$ws.Range("I86").Value = 2141.9167   # was 2300.0715
$ws.Range("J86").Value = 3400   # was 3300
$ws.Range("K86").Value = 2141.9167   # was 2300.0715
$ws.Range("L86").Value = 3400   # was 3300
$ws.Range("M86").Value = -1018.9167   # was -1177.0715
$ws.Range("N86").Value = -5646   # was -5546
$ws.Range("H89").Value = 2393.5334   # was 2425.0625
$ws.Range("I89").Value = 2141.9167   # was 2300.0715
$ws.Range("J89").Value = 3400   # was 3300
$ws.Range("K89").Value = 10709.5835   # was 11500.3575
$ws.Range("L89").Value = 17000   # was 16500
$ws.Range("M89").Value = -5093.583500000001   # was -5884.3575
$ws.Range("N89").Value = -28232   # was -27732

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3924.2163   # was 4362.3076
$ws.Range("I31").Value = 1577.0741   # was 1775.6086
$ws.Range("J31").Value = 5272.5747   # was 5778.8335
$ws.Range("K31").Value = 1577.0741   # was 1775.6086
$ws.Range("L31").Value = 5272.5747   # was 5778.8335
$ws.Range("M31").Value = -1282.0741   # was -1480.6086
$ws.Range("N31").Value = -5862.5747   # was -6368.8335
$ws.Range("H34").Value = 3924.2163   # was 4362.3076
$ws.Range("I34").Value = 1577.0741   # was 1775.6086
$ws.Range("J34").Value = 5272.5747   # was 5778.8335
$ws.Range("K34").Value = 1577.0741   # was 1775.6086
$ws.Range("L34").Value = 5272.5747   # was 5778.8335
$ws.Range("M34").Value = -1375.0741   # was -1573.6086
$ws.Range("N34").Value = -5676.5747   # was -6182.8335
$ws.Range("H50").Value = 24000   # was 19724
$ws.Range("J50").Value = 24000   # was 19724
$ws.Range("L50").Value = 24000   # was 19724
$ws.Range("N50").Value = -25250   # was -20974
$ws.Range("H59").Value = 24159   # was 25233
$ws.Range("J59").Value = 26448.75   # was 27279.6
$ws.Range("L59").Value = 26448.75   # was 27279.6
$ws.Range("N59").Value = -28738.75   # was -29569.6
$ws.Range("H60").Value = 19000   # was 18750
$ws.Range("J60").Value = 19000   # was 18750
$ws.Range("L60").Value = 19000   # was 18750
$ws.Range("N60").Value = -20022   # was -19772
$ws.Range("H68").Value = 32000   # was 0
$ws.Range("J68").Value = 32000   # was 0
$ws.Range("L68").Value = 32000   # was 0
$ws.Range("N68").Value = -33498   # new cell
$ws.Range("H71").Value = 32000   # was 0
$ws.Range("J71").Value = 32000   # was 0
$ws.Range("L71").Value = 96000   # was 0
$ws.Range("N71").Value = -103488   # new cell
$ws.Range("H74").Value = 28865.834   # was 28899.5
$ws.Range("J74").Value = 28865.834   # was 28899.5
$ws.Range("L74").Value = 28865.834   # was 28899.5
$ws.Range("N74").Value = -30613.834   # was -30647.5
$ws.Range("H77").Value = 28865.834   # was 28899.5
$ws.Range("J77").Value = 28865.834   # was 28899.5
$ws.Range("L77").Value = 86597.50199999999   # was 86698.5
$ws.Range("N77").Value = -95333.50199999999   # was -95434.5
$ws.Range("H99").Value = 9629772   # was 11380422
$ws.Range("I99").Value = 13417.714   # was 15504
$ws.Range("J99").Value = 20848852   # was 25018322
$ws.Range("K99").Value = 13417.714   # was 15504
$ws.Range("L99").Value = 20848852   # was 25018322
$ws.Range("M99").Value = -11919.714   # was -14006
$ws.Range("N99").Value = -20851848   # was -25021318
$ws.Range("H126").Value = 9629772   # was 11380422
$ws.Range("I126").Value = 13417.714   # was 15504
$ws.Range("J126").Value = 20848852   # was 25018322
$ws.Range("K126").Value = 40253.142   # was 46512
$ws.Range("L126").Value = 62546556   # was 75054966
$ws.Range("M126").Value = -37783.142   # was -44042
$ws.Range("N126").Value = -62551496   # was -75059906
$ws.Range("H132").Value = 5183.5454   # was 5401
$ws.Range("I132").Value = 4198.4   # was 5101
$ws.Range("J132").Value = 6004.5   # was 5572.4287
$ws.Range("K132").Value = 12595.2   # was 15303
$ws.Range("L132").Value = 18013.5   # was 16717.2861
$ws.Range("M132").Value = -10065.2   # was -12773
$ws.Range("N132").Value = -23073.5   # was -21777.2861

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5083.3335   # was 5928.5713
$ws.Range("I56").Value = 5083.3335   # was 5928.5713
$ws.Range("K56").Value = 5083.3335   # was 5928.5713
$ws.Range("M56").Value = -4553.3335   # was -5398.5713
$ws.Range("H80").Value = 6333.3335   # was 6250
$ws.Range("I80").Value = 3000   # was 2500
$ws.Range("K80").Value = 9000   # was 7500
$ws.Range("M80").Value = -8064   # was -6564
$ws.Range("H83").Value = 6333.3335   # was 6250
$ws.Range("I83").Value = 3000   # was 2500
$ws.Range("K83").Value = 27000   # was 22500
$ws.Range("M83").Value = -22320   # was -17820
$ws.Range("H113").Value = 151965.78   # was 1017400.3
$ws.Range("I113").Value = 432.1154   # was 1087385.8
$ws.Range("J113").Value = 714805.1   # was 769759.4399999999
$ws.Range("K113").Value = 1296.3462   # was 3262157.4
$ws.Range("L113").Value = 2144415.3   # was 2309278.32
$ws.Range("M113").Value = 873.6538   # was -3259987.4
$ws.Range("N113").Value = -2148755.3   # was -2313618.32
$ws.Range("H122").Value = 574.9286   # was 574.2143
$ws.Range("I122").Value = 464.4   # was 463.4
$ws.Range("K122").Value = 4179.599999999999   # was 4170.599999999999
$ws.Range("M122").Value = -1729.599999999999   # was -1720.599999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 52632720   # was 45455550
$ws.Range("I113").Value = 83334270   # was 66667590
$ws.Range("J113").Value = 1486.1428   # was 1184.2858
$ws.Range("K113").Value = 83334270   # was 66667590
$ws.Range("L113").Value = 1486.1428   # was 1184.2858
$ws.Range("M113").Value = -83332100   # was -66665420
$ws.Range("N113").Value = -5826.1428   # was -5524.2858

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4116.2915   # was 3831.173
$ws.Range("I136").Value = 2178.9429   # was 2018.4736
$ws.Range("J136").Value = 9332.23   # was 8751.357
$ws.Range("K136").Value = 6536.8287   # was 6055.4208
$ws.Range("L136").Value = 27996.69   # was 26254.071
$ws.Range("M136").Value = -3986.8287   # was -3505.4208
$ws.Range("N136").Value = -33096.69   # was -31354.071

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1332.2858   # was 1641.8334
$ws.Range("I126").Value = 990   # was 1283.1666
$ws.Range("J126").Value = 1788.6666   # was 2000.5
$ws.Range("K126").Value = 2970   # was 3849.4998
$ws.Range("L126").Value = 5365.9998   # was 6001.5
$ws.Range("M126").Value = -500   # was -1379.4998
$ws.Range("N126").Value = -10305.9998   # was -10941.5
$ws.Range("H136").Value = 2650.861   # was 2782.4119
$ws.Range("I136").Value = 2612.8147   # was 2788.68
$ws.Range("K136").Value = 7838.4441   # was 8366.039999999999
$ws.Range("M136").Value = -5288.4441   # was -5816.039999999999
